$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.258.39"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.909.20"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "307.86"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").Value = "0.3819"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").Value = "0.07314"
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("D10").Value = "21.63"
$ws.Range("E10").Value = "  +2.36%  "
$ws.Range("D11").Value = "0.9066"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "0.08232"
$ws.Range("E12").Value = "  -3.70%  "
$ws.Range("D13").Value = "96.47"
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("D14").Value = "5.377"
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").Value = "1.728.02"
$ws.Range("E15").Value = "  -9.23%  "
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "0.000008692"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "14.77"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "27.292.43"
$ws.Range("D21").Value = "5.133"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("D23").Value = "6.499"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("D24").Value = "2.346"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "1.740"
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").Value = "117.00"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").Value = "4.858"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").Value = "4.881"
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").Value = "0.09239"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").Value = "0.8272"
$ws.Range("E32").Value = "  +2.38%  "
$ws.Range("D33").Value = "0.05082"
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("D35").Value = "2.987"
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("D36").Value = "3.361"
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("D37").Value = "2.727"
$ws.Range("E37").Value = "  +4.03%  "
$ws.Range("D38").Value = "0.5761"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").Value = "0.02005"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("D41").Value = "9.064"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Value = "6.626"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").Value = "117.14"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").Value = "0.4939"
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("D46").Value = "10.19"
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("D49").Value = "38.70"
$ws.Range("E49").Value = "  +3.23%  "
$ws.Range("D50").Value = "64.54"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").Value = "0.06056"
$ws.Range("E51").Value = "  +2.10%  "
